$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.719.25"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").Value = "1.695.20"
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.40%  "

$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3957"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4068"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.493"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.004"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08901"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.52%  "

$ws.Range("E13").Value = "  -0.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.098"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001323"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("D17").Value = "1.697.64"
$ws.Range("E17").Value = "  +0.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07039"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.66%  "

$ws.Range("E20").Value = "  +0.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.013"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.006"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.51%  "

$ws.Range("D24").Value = "24.670.22"
$ws.Range("E24").Value = "  +0.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.219"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.362"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "136.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.193"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.559"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08636"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.056"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.070"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2736"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.48%  "

$ws.Range("B37").Value = "WEMIXTOKEN"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.888"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.13%  "

$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09230"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02727"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.472"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7684"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.605"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7167"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.224"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.48%  "

$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "140.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.322"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "90.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07982"
$ws.Range("D51").Style = "Normal"
